$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "2025-05-23 14:00"
$ws.Range("B3").Value = "U07EVCPHEMP"
$ws.Range("C3").Value = "Gaius Omonale"
$ws.Range("D3").Value = "Yankee"
$ws.Range("E3").Value = 50000
$ws.Range("F3").Value = 7839920123
$ws.Range("G3").Value = "Gaius Omonale"
$ws.Range("H3").Value = "Opay"
